# Appends the daily-generated portfolio rows (Colab export, 2025-06-20 16:00)
# to each of the three sheets, matching the committed OOXML diff.
$wb = $excel.ActiveWorkbook

# --- 大智投资组合收益: append 11 new rows starting at row 105 ---
$ws = $wb.Worksheets.Item('大智投资组合收益')
$ws.Range("A105").Value = '大智 (稳健智远)'
$ws.Range("B105").Value = '''000333'
$ws.Range("B105").Style = "Normal"
$ws.Range("C105").Value = '美的集团'
$ws.Range("D105").Value = 3.08
$ws.Range("E105").Value = 42.89719527444501
$ws.Range("F105").Value = 72.38
$ws.Range("G105").Value = 3104.89899396433
$ws.Range("H105").Value = 100876.2403810837
$ws.Range("I105").Value = '''202506201600'
$ws.Range("I105").Style = "Normal"

$ws.Range("A106").Value = '大智 (稳健智远)'
$ws.Range("B106").Value = '''510050'
$ws.Range("B106").Style = "Normal"
$ws.Range("C106").Value = '上证50ETF'
$ws.Range("D106").Value = 5.05
$ws.Range("E106").Value = 1852.638869852698
$ws.Range("F106").Value = 2.75
$ws.Range("G106").Value = 5094.75689209492
$ws.Range("H106").Value = 100876.2403810837
$ws.Range("I106").Value = '''202506201600'
$ws.Range("I106").Style = "Normal"

$ws.Range("A107").Value = '大智 (稳健智远)'
$ws.Range("B107").Value = '''510300'
$ws.Range("B107").Style = "Normal"
$ws.Range("C107").Value = '沪深300ETF'
$ws.Range("D107").Value = 4.91
$ws.Range("E107").Value = 1276.881426590205
$ws.Range("F107").Value = 3.88
$ws.Range("G107").Value = 4954.299935169995
$ws.Range("H107").Value = 100876.2403810837
$ws.Range("I107").Value = '''202506201600'
$ws.Range("I107").Style = "Normal"

$ws.Range("A108").Value = '大智 (稳健智远)'
$ws.Range("B108").Value = '''518880'
$ws.Range("B108").Style = "Normal"
$ws.Range("C108").Value = '黄金ETF'
$ws.Range("D108").Value = 4.97
$ws.Range("E108").Value = 673.9096418114974
$ws.Range("F108").Value = 7.44
$ws.Range("G108").Value = 5013.887735077541
$ws.Range("H108").Value = 100876.2403810837
$ws.Range("I108").Value = '''202506201600'
$ws.Range("I108").Style = "Normal"

$ws.Range("A109").Value = '大智 (稳健智远)'
$ws.Range("B109").Value = '''600085'
$ws.Range("B109").Style = "Normal"
$ws.Range("C109").Value = '同仁堂'
$ws.Range("D109").Value = 2
$ws.Range("E109").Value = 56.06334956913254
$ws.Range("F109").Value = 35.99
$ws.Range("G109").Value = 2017.71995099308
$ws.Range("H109").Value = 100876.2403810837
$ws.Range("I109").Value = '''202506201600'
$ws.Range("I109").Style = "Normal"

$ws.Range("A110").Value = '大智 (稳健智远)'
$ws.Range("B110").Value = '''600900'
$ws.Range("B110").Style = "Normal"
$ws.Range("C110").Value = '长江电力'
$ws.Range("D110").Value = 20.16
$ws.Range("E110").Value = 669.0422707938175
$ws.Range("F110").Value = 30.4
$ws.Range("G110").Value = 20338.88503213205
$ws.Range("H110").Value = 100876.2403810837
$ws.Range("I110").Value = '''202506201600'
$ws.Range("I110").Style = "Normal"

$ws.Range("A111").Value = '大智 (稳健智远)'
$ws.Range("B111").Value = '''600989'
$ws.Range("B111").Style = "Normal"
$ws.Range("C111").Value = '宝丰能源'
$ws.Range("D111").Value = 5.05
$ws.Range("E111").Value = 308.5861230826723
$ws.Range("F111").Value = 16.5
$ws.Range("G111").Value = 5091.671030864093
$ws.Range("H111").Value = 100876.2403810837
$ws.Range("I111").Value = '''202506201600'
$ws.Range("I111").Style = "Normal"

$ws.Range("A112").Value = '大智 (稳健智远)'
$ws.Range("B112").Value = '''601899'
$ws.Range("B112").Style = "Normal"
$ws.Range("C112").Value = 'XD紫金矿'
$ws.Range("D112").Value = 9.880000000000001
$ws.Range("E112").Value = 541.4194359293219
$ws.Range("F112").Value = 18.41
$ws.Range("G112").Value = 9967.531815458817
$ws.Range("H112").Value = 100876.2403810837
$ws.Range("I112").Value = '''202506201600'
$ws.Range("I112").Style = "Normal"

$ws.Range("A113").Value = '大智 (稳健智远)'
$ws.Range("B113").Value = '''HK02899'
$ws.Range("B113").Style = "Normal"
$ws.Range("C113").Value = '紫金矿业'
$ws.Range("D113").Value = 9.91
$ws.Range("E113").Value = 521.4694874201556
$ws.Range("F113").Value = 19.18
$ws.Range("G113").Value = 10001.78476871858
$ws.Range("H113").Value = 100876.2403810837
$ws.Range("I113").Value = '''202506201600'
$ws.Range("I113").Style = "Normal"

$ws.Range("A114").Value = '大智 (稳健智远)'
$ws.Range("B114").Value = '''HK06881'
$ws.Range("B114").Style = "Normal"
$ws.Range("C114").Value = '中国银河'
$ws.Range("D114").Value = 4.68
$ws.Range("E114").Value = 600.7968033130801
$ws.Range("F114").Value = 7.86
$ws.Range("G114").Value = 4722.26287404081
$ws.Range("H114").Value = 100876.2403810837
$ws.Range("I114").Value = '''202506201600'
$ws.Range("I114").Style = "Normal"

$ws.Range("A115").Value = '大智 (稳健智远)'
$ws.Range("B115").Value = '''100000'
$ws.Range("B115").Style = "Normal"
$ws.Range("C115").Value = '现金'
$ws.Range("D115").Value = 30.3
$ws.Range("E115").Value = 30568.54135256952
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 30568.54135256952
$ws.Range("H115").Value = 100876.2403810837
$ws.Range("I115").Value = '''202506201600'
$ws.Range("I115").Style = "Normal"


# --- 大成投资组合收益: append 7 new rows starting at row 72 ---
$ws = $wb.Worksheets.Item('大成投资组合收益')
$ws.Range("A72").Value = '大成 (锐进先锋)'
$ws.Range("B72").Value = '''000725'
$ws.Range("B72").Style = "Normal"
$ws.Range("C72").Value = '京东方A'
$ws.Range("D72").Value = 5.32
$ws.Range("E72").Value = 1264.047242468624
$ws.Range("F72").Value = 3.93
$ws.Range("G72").Value = 4967.705662901692
$ws.Range("H72").Value = 93344.42768391284
$ws.Range("I72").Value = '''202506201600'
$ws.Range("I72").Style = "Normal"

$ws.Range("A73").Value = '大成 (锐进先锋)'
$ws.Range("B73").Value = '''159781'
$ws.Range("B73").Style = "Normal"
$ws.Range("C73").Value = '科创创业ETF'
$ws.Range("D73").Value = 5.27
$ws.Range("E73").Value = 9277.629760760277
$ws.Range("F73").Value = 0.53
$ws.Range("G73").Value = 4917.143773202947
$ws.Range("H73").Value = 93344.42768391284
$ws.Range("I73").Value = '''202506201600'
$ws.Range("I73").Style = "Normal"

$ws.Range("A74").Value = '大成 (锐进先锋)'
$ws.Range("B74").Value = '''513100'
$ws.Range("B74").Style = "Normal"
$ws.Range("C74").Value = '纳指ETF'
$ws.Range("D74").Value = 5.23
$ws.Range("E74").Value = 3131.938709046463
$ws.Range("F74").Value = 1.56
$ws.Range("G74").Value = 4885.824386112482
$ws.Range("H74").Value = 93344.42768391284
$ws.Range("I74").Value = '''202506201600'
$ws.Range("I74").Style = "Normal"

$ws.Range("A75").Value = '大成 (锐进先锋)'
$ws.Range("B75").Value = '''513290'
$ws.Range("B75").Style = "Normal"
$ws.Range("C75").Value = '纳指生物科技ETF'
$ws.Range("D75").Value = 1.03
$ws.Range("E75").Value = 870.2909333102562
$ws.Range("F75").Value = 1.1
$ws.Range("G75").Value = 957.3200266412819
$ws.Range("H75").Value = 93344.42768391284
$ws.Range("I75").Value = '''202506201600'
$ws.Range("I75").Style = "Normal"

$ws.Range("A76").Value = '大成 (锐进先锋)'
$ws.Range("B76").Value = '''603119'
$ws.Range("B76").Style = "Normal"
$ws.Range("C76").Value = '浙江荣泰'
$ws.Range("D76").Value = 43.72
$ws.Range("E76").Value = 1069.978093782073
$ws.Range("F76").Value = 38.14
$ws.Range("G76").Value = 40808.96449684827
$ws.Range("H76").Value = 93344.42768391284
$ws.Range("I76").Value = '''202506201600'
$ws.Range("I76").Style = "Normal"

$ws.Range("A77").Value = '大成 (锐进先锋)'
$ws.Range("B77").Value = '''688290'
$ws.Range("B77").Style = "Normal"
$ws.Range("C77").Value = '景业智能'
$ws.Range("D77").Value = 7.83
$ws.Range("E77").Value = 147.121987895035
$ws.Range("F77").Value = 49.65
$ws.Range("G77").Value = 7304.606698988488
$ws.Range("H77").Value = 93344.42768391284
$ws.Range("I77").Value = '''202506201600'
$ws.Range("I77").Style = "Normal"

$ws.Range("A78").Value = '大成 (锐进先锋)'
$ws.Range("B78").Value = '''100000'
$ws.Range("B78").Style = "Normal"
$ws.Range("C78").Value = '现金'
$ws.Range("D78").Value = 31.61
$ws.Range("E78").Value = 29502.86263921768
$ws.Range("F78").Value = 1
$ws.Range("G78").Value = 29502.86263921768
$ws.Range("H78").Value = 93344.42768391284
$ws.Range("I78").Value = '''202506201600'
$ws.Range("I78").Style = "Normal"


# --- 我的投资组合收益: append 15 new rows starting at row 151 ---
$ws = $wb.Worksheets.Item('我的投资组合收益')
$ws.Range("A151").Value = '范式进化投资组合'
$ws.Range("B151").Value = '''000333'
$ws.Range("B151").Style = "Normal"
$ws.Range("C151").Value = '美的集团'
$ws.Range("D151").Value = 1.02
$ws.Range("E151").Value = 14.02515882310653
$ws.Range("F151").Value = 72.38
$ws.Range("G151").Value = 1015.140995616451
$ws.Range("H151").Value = 99440.34091078889
$ws.Range("I151").Value = '''202506201600'
$ws.Range("I151").Style = "Normal"

$ws.Range("A152").Value = '范式进化投资组合'
$ws.Range("B152").Value = '''000725'
$ws.Range("B152").Style = "Normal"
$ws.Range("C152").Value = '京东方A'
$ws.Range("D152").Value = 5.08
$ws.Range("E152").Value = 1284.618017653691
$ws.Range("F152").Value = 3.93
$ws.Range("G152").Value = 5048.548809379005
$ws.Range("H152").Value = 99440.34091078889
$ws.Range("I152").Value = '''202506201600'
$ws.Range("I152").Style = "Normal"

$ws.Range("A153").Value = '范式进化投资组合'
$ws.Range("B153").Value = '''159781'
$ws.Range("B153").Style = "Normal"
$ws.Range("C153").Value = '科创创业ETF'
$ws.Range("D153").Value = 5.03
$ws.Range("E153").Value = 9428.611488061997
$ws.Range("F153").Value = 0.53
$ws.Range("G153").Value = 4997.164088672858
$ws.Range("H153").Value = 99440.34091078889
$ws.Range("I153").Value = '''202506201600'
$ws.Range("I153").Style = "Normal"

$ws.Range("A154").Value = '范式进化投资组合'
$ws.Range("B154").Value = '''510050'
$ws.Range("B154").Style = "Normal"
$ws.Range("C154").Value = '上证50ETF'
$ws.Range("D154").Value = 5.03
$ws.Range("E154").Value = 1817.150577699221
$ws.Range("F154").Value = 2.75
$ws.Range("G154").Value = 4997.164088672857
$ws.Range("H154").Value = 99440.34091078889
$ws.Range("I154").Value = '''202506201600'
$ws.Range("I154").Style = "Normal"

$ws.Range("A155").Value = '范式进化投资组合'
$ws.Range("B155").Value = '''510300'
$ws.Range("B155").Style = "Normal"
$ws.Range("C155").Value = '沪深300ETF'
$ws.Range("D155").Value = 4.89
$ws.Range("E155").Value = 1252.422077361618
$ws.Range("F155").Value = 3.88
$ws.Range("G155").Value = 4859.397660163077
$ws.Range("H155").Value = 99440.34091078889
$ws.Range("I155").Value = '''202506201600'
$ws.Range("I155").Style = "Normal"

$ws.Range("A156").Value = '范式进化投资组合'
$ws.Range("B156").Value = '''513100'
$ws.Range("B156").Style = "Normal"
$ws.Range("C156").Value = '纳指ETF'
$ws.Range("D156").Value = 1
$ws.Range("E156").Value = 636.5814125697908
$ws.Range("F156").Value = 1.56
$ws.Range("G156").Value = 993.0670036088737
$ws.Range("H156").Value = 99440.34091078889
$ws.Range("I156").Value = '''202506201600'
$ws.Range("I156").Style = "Normal"

$ws.Range("A157").Value = '范式进化投资组合'
$ws.Range("B157").Value = '''513290'
$ws.Range("B157").Style = "Normal"
$ws.Range("C157").Value = '纳指生物科技ETF'
$ws.Range("D157").Value = 0.98
$ws.Range("E157").Value = 884.4538210040457
$ws.Range("F157").Value = 1.1
$ws.Range("G157").Value = 972.8992031044504
$ws.Range("H157").Value = 99440.34091078889
$ws.Range("I157").Value = '''202506201600'
$ws.Range("I157").Style = "Normal"

$ws.Range("A158").Value = '范式进化投资组合'
$ws.Range("B158").Value = '''518880'
$ws.Range("B158").Style = "Normal"
$ws.Range("C158").Value = '黄金ETF'
$ws.Range("D158").Value = 0.99
$ws.Range("E158").Value = 132.2001081659486
$ws.Range("F158").Value = 7.44
$ws.Range("G158").Value = 983.5688047546575
$ws.Range("H158").Value = 99440.34091078889
$ws.Range("I158").Value = '''202506201600'
$ws.Range("I158").Style = "Normal"

$ws.Range("A159").Value = '范式进化投资组合'
$ws.Range("B159").Value = '''600085'
$ws.Range("B159").Style = "Normal"
$ws.Range("C159").Value = '同仁堂'
$ws.Range("D159").Value = 1
$ws.Range("E159").Value = 27.49471300507762
$ws.Range("F159").Value = 35.99
$ws.Range("G159").Value = 989.5347210527436
$ws.Range("H159").Value = 99440.34091078889
$ws.Range("I159").Value = '''202506201600'
$ws.Range("I159").Style = "Normal"

$ws.Range("A160").Value = '范式进化投资组合'
$ws.Range("B160").Value = '''600900'
$ws.Range("B160").Style = "Normal"
$ws.Range("C160").Value = '长江电力'
$ws.Range("D160").Value = 1
$ws.Range("E160").Value = 32.81132034584936
$ws.Range("F160").Value = 30.4
$ws.Range("G160").Value = 997.4641385138204
$ws.Range("H160").Value = 99440.34091078889
$ws.Range("I160").Value = '''202506201600'
$ws.Range("I160").Style = "Normal"

$ws.Range("A161").Value = '范式进化投资组合'
$ws.Range("B161").Value = '''600989'
$ws.Range("B161").Style = "Normal"
$ws.Range("C161").Value = '宝丰能源'
$ws.Range("D161").Value = 5.02
$ws.Range("E161").Value = 302.6749902285196
$ws.Range("F161").Value = 16.5
$ws.Range("G161").Value = 4994.137338770573
$ws.Range("H161").Value = 99440.34091078889
$ws.Range("I161").Value = '''202506201600'
$ws.Range("I161").Style = "Normal"

$ws.Range("A162").Value = '范式进化投资组合'
$ws.Range("B162").Value = '''601899'
$ws.Range("B162").Style = "Normal"
$ws.Range("C162").Value = 'XD紫金矿'
$ws.Range("D162").Value = 9.83
$ws.Range("E162").Value = 531.0482559694855
$ws.Range("F162").Value = 18.41
$ws.Range("G162").Value = 9776.598392398228
$ws.Range("H162").Value = 99440.34091078889
$ws.Range("I162").Value = '''202506201600'
$ws.Range("I162").Style = "Normal"

$ws.Range("A163").Value = '范式进化投资组合'
$ws.Range("B163").Value = '''603119'
$ws.Range("B163").Style = "Normal"
$ws.Range("C163").Value = '浙江荣泰'
$ws.Range("D163").Value = 0.93
$ws.Range("E163").Value = 24.16423640557475
$ws.Range("F163").Value = 38.14
$ws.Range("G163").Value = 921.623976508621
$ws.Range("H163").Value = 99440.34091078889
$ws.Range("I163").Value = '''202506201600'
$ws.Range("I163").Style = "Normal"

$ws.Range("A164").Value = '范式进化投资组合'
$ws.Range("B164").Value = '''HK06881'
$ws.Range("B164").Style = "Normal"
$ws.Range("C164").Value = '中国银河'
$ws.Range("D164").Value = 0.93
$ws.Range("E164").Value = 117.8576436007749
$ws.Range("F164").Value = 7.86
$ws.Range("G164").Value = 926.3610787020908
$ws.Range("H164").Value = 99440.34091078889
$ws.Range("I164").Value = '''202506201600'
$ws.Range("I164").Style = "Normal"

$ws.Range("A165").Value = '范式进化投资组合'
$ws.Range("B165").Value = '''100000'
$ws.Range("B165").Style = "Normal"
$ws.Range("C165").Value = '现金'
$ws.Range("D165").Value = 57.29
$ws.Range("E165").Value = 56967.67061087058
$ws.Range("F165").Value = 1
$ws.Range("G165").Value = 56967.67061087058
$ws.Range("H165").Value = 99440.34091078889
$ws.Range("I165").Value = '''202506201600'
$ws.Range("I165").Style = "Normal"

